$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Wafungwa na peremende - manukuu:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Prisoners and candies - subtitles:", 2)

$d.Content.Find.Execute(
    "**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino", $true, $false, $false, $false, $false,
    $true, 1, $false, "**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino", 2)

$d.Content.Find.Execute(
    "[Muziki]", $true, $false, $false, $false, $false,
    $true, 1, $false, "[Music]", 2)

$d.Save()
